$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E (particip) and F (taxa_sucesso) for rows 2-7 are being
# converted from fractional values (0-1) to percentage-scale values (0-100).
for ($row = 2; $row -le 7; $row++) {
    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)
    $eCell.Value = $eCell.Value() * 100
    $fCell.Value = $fCell.Value() * 100
}
